# DDT bug finding is done.
# Remove the old header row (FirstName/LastName/.../Password) so the data
# that used to be in row 2 (nayeem01 ...) becomes row 1, shifting every
# subsequent row up by one and dropping the former last row (nayeem06/user35).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(1).Delete()

# Rename the UserName column values from user30..user35 to test01..test06.
$ws.Range("I1").Value = "test01"
$ws.Range("I2").Value = "test02"
$ws.Range("I3").Value = "test03"
$ws.Range("I4").Value = "test04"
$ws.Range("I5").Value = "test05"
$ws.Range("I6").Value = "test06"

# Match the saved selection state.
$ws.Range("H9").Select()
